$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 279.72223
$ws.Range("I2").Value = 291.52942
$ws.Range("K2").Value = 291.52942
$ws.Range("M2").Value = -178.52942
$ws.Range("H17").Value = 1129.1351
$ws.Range("J17").Value = 1129.1351
$ws.Range("L17").Value = 3387.4053
$ws.Range("N17").Value = -3723.4053
$ws.Range("H76").Value = 3266.5186
$ws.Range("I76").Value = 3087.0476
$ws.Range("K76").Value = 3087.0476
$ws.Range("M76").Value = -2772.0476
$ws.Range("H79").Value = 3266.5186
$ws.Range("I79").Value = 3087.0476
$ws.Range("K79").Value = 3087.0476
$ws.Range("M79").Value = -1995.0476
$ws.Range("H92").Value = 4775.6313
$ws.Range("I92").Value = 2107.6155
$ws.Range("K92").Value = 2107.6155
$ws.Range("M92").Value = -859.6154999999999
$ws.Range("H100").Value = 6320.923
$ws.Range("J100").Value = 7596
$ws.Range("L100").Value = 7596
$ws.Range("N100").Value = -8678
$ws.Range("H112").Value = 3019.6
$ws.Range("J112").Value = 5249
$ws.Range("L112").Value = 15747
$ws.Range("N112").Value = -17963
$ws.Range("H138").Value = 2234.3804
$ws.Range("J138").Value = 2660.0588
$ws.Range("L138").Value = 7980.176399999999
$ws.Range("N138").Value = -18260.1764
$ws.Range("H141").Value = 5854.033
$ws.Range("I141").Value = 3453.7896
$ws.Range("J141").Value = 9999.909
$ws.Range("K141").Value = 10361.3688
$ws.Range("L141").Value = 29999.727
$ws.Range("M141").Value = -5181.3688
$ws.Range("N141").Value = -40359.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 131.5
$ws.Range("I4").Value = 160.5
$ws.Range("K4").Value = 160.5
$ws.Range("M4").Value = -44.5
$ws.Range("H5").Value = 179.44444
$ws.Range("I5").Value = 69.333336
$ws.Range("J5").Value = 399.66666
$ws.Range("K5").Value = 69.333336
$ws.Range("L5").Value = 399.66666
$ws.Range("M5").Value = 42.666664
$ws.Range("N5").Value = -623.66666
$ws.Range("H32").Value = 2369.0784
$ws.Range("I32").Value = 2415.38
$ws.Range("K32").Value = 2415.38
$ws.Range("M32").Value = -2128.38
$ws.Range("H61").Value = 4828
$ws.Range("I61").Value = 4369.4
$ws.Range("K61").Value = 4369.4
$ws.Range("M61").Value = -4157.4
$ws.Range("H102").Value = 4014.9285
$ws.Range("I102").Value = 2621
$ws.Range("K102").Value = 2621
$ws.Range("M102").Value = -999
$ws.Range("H110").Value = 951.8
$ws.Range("I110").Value = 963.1111
$ws.Range("J110").Value = 850
$ws.Range("K110").Value = 963.1111
$ws.Range("L110").Value = 850
$ws.Range("M110").Value = 1081.8889
$ws.Range("N110").Value = -4940
$ws.Range("H136").Value = 4828
$ws.Range("I136").Value = 4369.4
$ws.Range("K136").Value = 13108.2
$ws.Range("M136").Value = -10558.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 179.44444
$ws.Range("I4").Value = 69.333336
$ws.Range("J4").Value = 399.66666
$ws.Range("K4").Value = 69.333336
$ws.Range("L4").Value = 399.66666
$ws.Range("M4").Value = 45.666664
$ws.Range("N4").Value = -629.66666
$ws.Range("H105").Value = 3151.077
$ws.Range("I105").Value = 2459.1428
$ws.Range("J105").Value = 3958.3333
$ws.Range("K105").Value = 2459.1428
$ws.Range("L105").Value = 3958.3333
$ws.Range("M105").Value = -712.1428000000001
$ws.Range("N105").Value = -7452.3333
$ws.Range("H134").Value = 3092.3157
$ws.Range("I134").Value = 3109.625
$ws.Range("K134").Value = 9328.875
$ws.Range("M134").Value = -6793.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 302.25
$ws.Range("I5").Value = 302.5
$ws.Range("K5").Value = 302.5
$ws.Range("M5").Value = -190.5
$ws.Range("H23").Value = 5250
$ws.Range("J23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("N23").Value = -5480
$ws.Range("H27").Value = 5250
$ws.Range("J27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("N27").Value = -5384
$ws.Range("H55").Value = 12997.714
$ws.Range("J55").Value = 12997.714
$ws.Range("L55").Value = 12997.714
$ws.Range("N55").Value = -13627.714
$ws.Range("H107").Value = 42862.543
$ws.Range("I107").Value = 125415.25
$ws.Range("J107").Value = 1586.1875
$ws.Range("K107").Value = 125415.25
$ws.Range("L107").Value = 1586.1875
$ws.Range("M107").Value = -123495.25
$ws.Range("N107").Value = -5426.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 34.52
$ws.Range("J2").Value = 67.2
$ws.Range("L2").Value = 403.2
$ws.Range("N2").Value = -629.2
$ws.Range("H121").Value = 56252.445
$ws.Range("I121").Value = 91483.63
$ws.Range("J121").Value = 889.1429000000001
$ws.Range("K121").Value = 274450.89
$ws.Range("L121").Value = 2667.4287
$ws.Range("M121").Value = -273140.89
$ws.Range("N121").Value = -5287.4287
$ws.Range("H129").Value = 2248.0908
$ws.Range("J129").Value = 4032.7
$ws.Range("L129").Value = 12098.1
$ws.Range("N129").Value = -22098.1
$ws.Range("H131").Value = 1782731.9
$ws.Range("J131").Value = 2780450
$ws.Range("L131").Value = 8341350
$ws.Range("N131").Value = -8351430
$ws.Range("H137").Value = 8870.429
$ws.Range("I137").Value = 2866.3333
$ws.Range("J137").Value = 13373.5
$ws.Range("K137").Value = 8598.999899999999
$ws.Range("L137").Value = 40120.5
$ws.Range("M137").Value = -3498.999899999999
$ws.Range("N137").Value = -50320.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 12525000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H12").Value = 1039874.4
$ws.Range("I12").Value = 1002713.56
$ws.Range("K12").Value = 1002713.56
$ws.Range("M12").Value = -1002573.56
$ws.Range("H14").Value = 63950
$ws.Range("I14").Value = 83600
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 83600
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = -83432
$ws.Range("N14").Value = -5336
$ws.Range("H52").Value = 515000
$ws.Range("J52").Value = 515000
$ws.Range("L52").Value = 515000
$ws.Range("N52").Value = -515518
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H122").Value = 9156.299999999999
$ws.Range("I122").Value = 9392.296
$ws.Range("J122").Value = 7032.3335
$ws.Range("K122").Value = 28176.888
$ws.Range("L122").Value = 21097.0005
$ws.Range("M122").Value = -25726.888
$ws.Range("N122").Value = -25997.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 2100
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H15").Value = 2100
$ws.Range("H17").Value = 4666.5
$ws.Range("I17").Value = 2666.6667
$ws.Range("J17").Value = 10666
$ws.Range("K17").Value = 2666.6667
$ws.Range("L17").Value = 10666
$ws.Range("M17").Value = -2496.6667
$ws.Range("N17").Value = -11006
$ws.Range("H22").Value = 2533.524
$ws.Range("I22").Value = 2261.5789
$ws.Range("J22").Value = 2758.1738
$ws.Range("K22").Value = 2261.5789
$ws.Range("L22").Value = 2758.1738
$ws.Range("M22").Value = -1966.5789
$ws.Range("N22").Value = -3348.1738
$ws.Range("H27").Value = 2533.524
$ws.Range("I27").Value = 2261.5789
$ws.Range("J27").Value = 2758.1738
$ws.Range("K27").Value = 2261.5789
$ws.Range("L27").Value = 2758.1738
$ws.Range("M27").Value = -2154.5789
$ws.Range("N27").Value = -2972.1738
$ws.Range("H68").Value = 7636.364
$ws.Range("I68").Value = 3666.6667
$ws.Range("J68").Value = 9125
$ws.Range("K68").Value = 3666.6667
$ws.Range("L68").Value = 9125
$ws.Range("M68").Value = -2917.6667
$ws.Range("N68").Value = -10623
$ws.Range("H71").Value = 7636.364
$ws.Range("I71").Value = 3666.6667
$ws.Range("J71").Value = 9125
$ws.Range("K71").Value = 18333.3335
$ws.Range("L71").Value = 45625
$ws.Range("M71").Value = -14589.3335
$ws.Range("N71").Value = -53113
$ws.Range("H122").Value = 3899
$ws.Range("I122").Value = 3899
$ws.Range("K122").Value = 11697
$ws.Range("M122").Value = -9247
$ws.Range("H132").Value = 3888.4211
$ws.Range("I132").Value = 4143.737
$ws.Range("J132").Value = 3633.1052
$ws.Range("K132").Value = 12431.211
$ws.Range("L132").Value = 10899.3156
$ws.Range("M132").Value = -9901.210999999999
$ws.Range("N132").Value = -15959.3156

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 9500
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("N17").Value = -10344
$ws.Range("H70").Value = 35471
$ws.Range("I70").Value = 19000
$ws.Range("J70").Value = 40961.332
$ws.Range("K70").Value = 19000
$ws.Range("L70").Value = 40961.332
$ws.Range("M70").Value = -18685
$ws.Range("N70").Value = -41591.332
$ws.Range("H73").Value = 35471
$ws.Range("I73").Value = 19000
$ws.Range("J73").Value = 40961.332
$ws.Range("K73").Value = 19000
$ws.Range("L73").Value = 40961.332
$ws.Range("M73").Value = -17908
$ws.Range("N73").Value = -43145.332
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H113").Value = 1251.7333
$ws.Range("I113").Value = 1135.3478
$ws.Range("J113").Value = 1634.1428
$ws.Range("K113").Value = 3406.0434
$ws.Range("L113").Value = 4902.428400000001
$ws.Range("M113").Value = -1236.0434
$ws.Range("N113").Value = -9242.428400000001
$ws.Range("H132").Value = 3520.8
$ws.Range("I132").Value = 3410.6667
$ws.Range("K132").Value = 10232.0001
$ws.Range("M132").Value = -7702.000100000001
